$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.083.79"
$ws.Range("E2").Value = "  +4.52%  "

# Row 3
$ws.Range("D3").Value = "2.238.68"
$ws.Range("E3").Value = "  +4.42%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "'251.12"
$ws.Range("E5").Value = "  +6.65%  "

# Row 6
$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = "  +2.60%  "

# Row 7
$ws.Range("D7").Value = "'75.06"
$ws.Range("E7").Value = "  +8.85%  "

# Row 8
$ws.Range("E8").Value = "  -0.15%  "

# Row 9
$ws.Range("D9").Value = "'0.599"
$ws.Range("E9").Value = "  +6.30%  "

# Row 10
$ws.Range("D10").Value = "'40.93"
$ws.Range("E10").Value = "  +7.04%  "

# Row 11
$ws.Range("D11").Value = "'0.0930"
$ws.Range("E11").Value = "  +4.15%  "

# Row 12
$ws.Range("D12").Value = "'6.91"
$ws.Range("E12").Value = "  +4.99%  "

# Row 13
$ws.Range("D13").Value = "'0.101"
$ws.Range("E13").Value = "  +1.64%  "

# Row 14
$ws.Range("D14").Value = "2.577.74"
$ws.Range("E14").Value = "  +4.59%  "

# Row 15
$ws.Range("D15").Value = "'14.61"
$ws.Range("E15").Value = "  +1.57%  "

# Row 16
$ws.Range("D16").Value = "2.246.01"
$ws.Range("E16").Value = "  +3.77%  "

# Row 17
$ws.Range("E17").Value = "  +2.23%  "

# Row 18
$ws.Range("D18").Value = "42.999.10"
$ws.Range("E18").Value = "  +4.79%  "

# Row 19
$ws.Range("E19").Value = "  +5.49%  "

# Row 20
$ws.Range("D20").Value = "'71.15"
$ws.Range("E20").Value = "  +2.95%  "

# Row 21
$ws.Range("D21").Value = "'5.99"
$ws.Range("E21").Value = "  +4.97%  "

# Row 22
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'9.74"
$ws.Range("E22").Value = "  +3.02%  "

# Row 23
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'229.73"
$ws.Range("E23").Value = "  +2.59%  "

# Row 24
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "'2.19"
$ws.Range("E24").Value = "  +17.01%  "

# Row 25
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
$ws.Range("D26").Value = "'10.78"
$ws.Range("E26").Value = "  +2.45%  "

# Row 27
$ws.Range("E27").Value = "  +2.37%  "

# Row 28
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = "  +5.76%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'39.44"
$ws.Range("E29").Value = "  +28.72%  "

# Row 30
$ws.Range("E30").Value = "  +4.07%  "

# Row 31
$ws.Range("D31").Value = "'171.58"
$ws.Range("E31").Value = "  +1.83%  "

# Row 32
$ws.Range("D32").Value = "'20.23"
$ws.Range("E32").Value = "  +3.49%  "

# Row 33
$ws.Range("D33").Value = "'0.0800"
$ws.Range("E33").Value = "  +6.81%  "

# Row 34
$ws.Range("D34").Value = "'5.25"
$ws.Range("E34").Value = "  +4.63%  "

# Row 35
$ws.Range("E35").Value = "  +2.39%  "

# Row 36
$ws.Range("E36").Value = "  +10.47%  "

# Row 37
$ws.Range("E37").Value = "  +10.28%  "

# Row 38
$ws.Range("D38").Value = "'0.0330"
$ws.Range("E38").Value = "  +18.46%  "

# Row 39
$ws.Range("D39").Value = "'12.99"
$ws.Range("E39").Value = "  +12.07%  "

# Row 40
$ws.Range("E40").Value = "  +4.30%  "

# Row 41
$ws.Range("D41").Value = "'0.205"
$ws.Range("E41").Value = "  +11.16%  "

# Row 42
$ws.Range("E42").Value = "  +3.51%  "

# Row 43
$ws.Range("D43").Value = "'59.42"
$ws.Range("E43").Value = "  +4.09%  "

# Row 44
$ws.Range("D44").Value = "'0.491"
$ws.Range("E44").Value = "  +33.94%  "

# Row 45
$ws.Range("D45").Value = "'8.64"
$ws.Range("E45").Value = "  +6.25%  "

# Row 46
$ws.Range("D46").Value = "'103.81"
$ws.Range("E46").Value = "  +7.75%  "

# Row 47
$ws.Range("E47").Value = "  +4.63%  "

# Row 48
$ws.Range("E48").Value = "  +13.53%  "

# Row 49
$ws.Range("D49").Value = "'1.10"
$ws.Range("E49").Value = "  +3.65%  "

# Row 50
$ws.Range("E50").Value = "  +4.98%  "

# Row 51
$ws.Range("D51").Value = "'2.68"
$ws.Range("E51").Value = "  +3.63%  "

